# "removed customers xcel file and copied over to vendors"
# Adds a second worksheet ("Sheet2" - the former "customers" sheet / shopping
# list) after the existing "Sheet1" (vendors) sheet, makes it the active
# sheet/tab, and restores Sheet1's prior selection.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 (Excel's default Add() would put
# it before the active sheet, so pass Before=$null, After=$ws1).
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Worksheets.Add() clones the active sheet's contents/formatting - wipe it
# back to blank before writing the new "shopping list" data.
$ws2.UsedRange.ClearContents() | Out-Null

$ws2.Range("A1").Value = "Name"
$ws2.Range("B1").Value = "Shopping List"
$ws2.Range("C1").Value = "Number"

$ws2.Range("A2").Value = "Bob"
$ws2.Range("B2").Value = "Bananas"
$ws2.Range("B3").Value = "Socks"
$ws2.Range("B4").Value = "T.V"

# Restore/assign the selections seen on each tab.
$ws1.Range("B13").Select() | Out-Null
$ws2.Range("E3").Select() | Out-Null
